function Find-ParaIndexByText($doc, $text) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $t = $doc.Paragraphs($i).Range.Text.Trim()
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Underline the title "Deliverable 2 Report"
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.Font.Underline = 1

# ---------------------------------------------------------------------------
# 2. Insert the "Team Members:" heading + the six team-member lines before
#    the existing "Team Assignments:" paragraph.
# ---------------------------------------------------------------------------
$taIndex = Find-ParaIndexByText $d "Team Assignments:"
$taPara = $d.Paragraphs($taIndex)
$taPara.Range.InsertParagraphBefore()

$newBlock = $d.Paragraphs($taIndex)
$newBlock.Range.Text = "Team Members:`rWai Fong – 11382065, kuanrya000`rShrunga Mallavalli – 11436985, malaval21`rLinh Nguyen – 11563329, linhnguyen14a2`rCary Ott – 11440278, CarlyOtt`rKimi Phan – 11466435, kphanswims15`rKayla Rhodes – 11373485, rhodeskl"

# Bold the "Team Members:" heading paragraph.
$teamMembersPara = $d.Paragraphs($taIndex)
$teamMembersPara.Range.Font.Bold = 1

# Bold the (now shifted) "Team Assignments:" paragraph.
$taIndex2 = Find-ParaIndexByText $d "Team Assignments:"
$taPara2 = $d.Paragraphs($taIndex2)
$taPara2.Range.Font.Bold = 1

# ---------------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark from the end of the document to sit right
#    after "Kimi Phan – 11466435, " (Word keeps only one "_GoBack" bookmark,
#    so re-adding it elsewhere relocates it and drops the old one).
# ---------------------------------------------------------------------------
$bmRange = $d.Content
$bmRange.Find.Execute("Kimi Phan – 11466435, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 4. Remove the stray "lastRenderedPageBreak" marker that used to sit before
#    the standalone "ExpressionCounter" heading.
# ---------------------------------------------------------------------------
$ecIndex = Find-ParaIndexByText $d "ExpressionCounter"
$ecPara = $d.Paragraphs($ecIndex)
$ecPara.Range.Text = "ExpressionCounter"

Write-Output "done"
